$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 17. This pushes the existing rows 17-27 down to 18-28
# (matching the diff's row-shift + new dimension A1:T28), and carries the
# date-format style from the cell above into the new D17.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new weekly record. Columns
# A, B, C, E, F, G, H, I, J, K, L repeat the same "template" values shared
# by every row in this block.
$ws.Cells.Item(17, 1).Value = 10
$ws.Cells.Item(17, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(17, 3).Value = "La Araucanía"
$ws.Cells.Item(17, 4).Value = 44658
$ws.Cells.Item(17, 5).Value = 9
$ws.Cells.Item(17, 6).Value = "Fruta"
$ws.Cells.Item(17, 7).Value = 100108
$ws.Cells.Item(17, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(17, 9).Value = 100108003
$ws.Cells.Item(17, 10).Value = "Maracuyá"
$ws.Cells.Item(17, 11).Value = "Sin especificar"
$ws.Cells.Item(17, 12).Value = "Primera"
$ws.Cells.Item(17, 13).Value = 30
$ws.Cells.Item(17, 14).Value = 28000
$ws.Cells.Item(17, 15).Value = 28000
$ws.Cells.Item(17, 16).Value = 28000
$ws.Cells.Item(17, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(17, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(17, 19).Value = 1556
$ws.Cells.Item(17, 20).Value = 18
